# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet to the latest scraped values, matching the upstream commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the text looks like
# a plain decimal number (e.g. "500.70"). Such values must be written with
# the cell pre-formatted as Text ("@"), otherwise Excel silently reinterprets
# them as floating point numbers (e.g. "500.70" -> 500.6999999999999) and the
# exact original text (including trailing zeros) would be lost.
$updates = @(
    @{ Cell = "D2"; Value = '55.775.27'; Numeric = $false }
    @{ Cell = "E2"; Value = '  -2.15%  '; Numeric = $false }
    @{ Cell = "D3"; Value = '2.967.26'; Numeric = $false }
    @{ Cell = "E3"; Value = '  -0.49%  '; Numeric = $false }
    @{ Cell = "D4"; Value = '0.999'; Numeric = $true }
    @{ Cell = "E4"; Value = '  -0.10%  '; Numeric = $false }
    @{ Cell = "D5"; Value = '500.70'; Numeric = $true }
    @{ Cell = "E5"; Value = '  -0.14%  '; Numeric = $false }
    @{ Cell = "D6"; Value = '136.30'; Numeric = $true }
    @{ Cell = "E6"; Value = '  -1.40%  '; Numeric = $false }
    @{ Cell = "E7"; Value = '  +0.03%  '; Numeric = $false }
    @{ Cell = "D8"; Value = '0.424'; Numeric = $true }
    @{ Cell = "E8"; Value = '  -1.43%  '; Numeric = $false }
    @{ Cell = "D9"; Value = '7.13'; Numeric = $true }
    @{ Cell = "E9"; Value = '  -2.54%  '; Numeric = $false }
    @{ Cell = "D10"; Value = '0.106'; Numeric = $true }
    @{ Cell = "E10"; Value = '  -1.76%  '; Numeric = $false }
    @{ Cell = "D11"; Value = '0.361'; Numeric = $true }
    @{ Cell = "E11"; Value = '  +0.57%  '; Numeric = $false }
    @{ Cell = "D12"; Value = '3.478.31'; Numeric = $false }
    @{ Cell = "E12"; Value = '  -0.24%  '; Numeric = $false }
    @{ Cell = "E13"; Value = '  -1.55%  '; Numeric = $false }
    @{ Cell = "D14"; Value = '25.79'; Numeric = $true }
    @{ Cell = "E14"; Value = '  -1.00%  '; Numeric = $false }
    @{ Cell = "D15"; Value = '0.0000159'; Numeric = $true }
    @{ Cell = "E15"; Value = '  -0.49%  '; Numeric = $false }
    @{ Cell = "D16"; Value = '55.769.11'; Numeric = $false }
    @{ Cell = "E16"; Value = '  -2.32%  '; Numeric = $false }
    @{ Cell = "D17"; Value = '2.961.96'; Numeric = $false }
    @{ Cell = "E17"; Value = '  -0.73%  '; Numeric = $false }
    @{ Cell = "D18"; Value = '5.96'; Numeric = $true }
    @{ Cell = "E18"; Value = '  -1.50%  '; Numeric = $false }
    @{ Cell = "D19"; Value = '12.79'; Numeric = $true }
    @{ Cell = "E19"; Value = '  +1.29%  '; Numeric = $false }
    @{ Cell = "D20"; Value = '7.94'; Numeric = $true }
    @{ Cell = "E20"; Value = '  +0.82%  '; Numeric = $false }
    @{ Cell = "D21"; Value = '326.63'; Numeric = $true }
    @{ Cell = "E21"; Value = '  +1.76%  '; Numeric = $false }
    @{ Cell = "E22"; Value = '  +0.01%  '; Numeric = $false }
    @{ Cell = "D23"; Value = '0.490'; Numeric = $true }
    @{ Cell = "E23"; Value = '  +0.16%  '; Numeric = $false }
    @{ Cell = "D24"; Value = '64.21'; Numeric = $true }
    @{ Cell = "E24"; Value = '  +0.66%  '; Numeric = $false }
    @{ Cell = "D25"; Value = '3.090.78'; Numeric = $false }
    @{ Cell = "E25"; Value = '  -0.31%  '; Numeric = $false }
    @{ Cell = "D26"; Value = '1.00'; Numeric = $true }
    @{ Cell = "E26"; Value = '  +0.19%  '; Numeric = $false }
    @{ Cell = "D27"; Value = '0.161'; Numeric = $true }
    @{ Cell = "E27"; Value = '  -2.01%  '; Numeric = $false }
    @{ Cell = "D28"; Value = '0.0₃0886'; Numeric = $false }
    @{ Cell = "E28"; Value = '  -1.32%  '; Numeric = $false }
    @{ Cell = "D29"; Value = '6.32'; Numeric = $true }
    @{ Cell = "E29"; Value = '  -3.44%  '; Numeric = $false }
    @{ Cell = "D30"; Value = '6.91'; Numeric = $true }
    @{ Cell = "E30"; Value = '  -2.14%  '; Numeric = $false }
    @{ Cell = "D31"; Value = '1.76'; Numeric = $true }
    @{ Cell = "E31"; Value = '  -0.76%  '; Numeric = $false }
    @{ Cell = "D32"; Value = '20.05'; Numeric = $true }
    @{ Cell = "E32"; Value = '  -0.68%  '; Numeric = $false }
    @{ Cell = "D33"; Value = '1.14'; Numeric = $true }
    @{ Cell = "E33"; Value = '  -2.21%  '; Numeric = $false }
    @{ Cell = "D34"; Value = '153.46'; Numeric = $true }
    @{ Cell = "E34"; Value = '  -1.25%  '; Numeric = $false }
    @{ Cell = "D35"; Value = '4.45'; Numeric = $true }
    @{ Cell = "E35"; Value = '  -2.81%  '; Numeric = $false }
    @{ Cell = "D36"; Value = '5.65'; Numeric = $true }
    @{ Cell = "E36"; Value = '  -2.34%  '; Numeric = $false }
    @{ Cell = "D37"; Value = '24.93'; Numeric = $true }
    @{ Cell = "E37"; Value = '  +3.60%  '; Numeric = $false }
    @{ Cell = "E38"; Value = '  -1.89%  '; Numeric = $false }
    @{ Cell = "D39"; Value = '0.0653'; Numeric = $true }
    @{ Cell = "E39"; Value = '  -2.10%  '; Numeric = $false }
    @{ Cell = "D40"; Value = '2.998.11'; Numeric = $false }
    @{ Cell = "E40"; Value = '  -0.52%  '; Numeric = $false }
    @{ Cell = "D41"; Value = '36.61'; Numeric = $true }
    @{ Cell = "D42"; Value = '0.999'; Numeric = $true }
    @{ Cell = "E42"; Value = '  -0.13%  '; Numeric = $false }
    @{ Cell = "D43"; Value = '3.74'; Numeric = $true }
    @{ Cell = "E43"; Value = '  -0.04%  '; Numeric = $false }
    @{ Cell = "D44"; Value = '0.645'; Numeric = $true }
    @{ Cell = "E44"; Value = '  +0.68%  '; Numeric = $false }
    @{ Cell = "D45"; Value = '2.146.31'; Numeric = $false }
    @{ Cell = "E45"; Value = '  -2.47%  '; Numeric = $false }
    @{ Cell = "D46"; Value = '1.34'; Numeric = $true }
    @{ Cell = "E46"; Value = '  -4.08%  '; Numeric = $false }
    @{ Cell = "D47"; Value = '5.79'; Numeric = $true }
    @{ Cell = "E47"; Value = '  -3.01%  '; Numeric = $false }
    @{ Cell = "D48"; Value = '0.914'; Numeric = $true }
    @{ Cell = "E48"; Value = '  -3.45%  '; Numeric = $false }
    @{ Cell = "D49"; Value = '0.0233'; Numeric = $true }
    @{ Cell = "E49"; Value = '  -0.94%  '; Numeric = $false }
    @{ Cell = "D50"; Value = '19.42'; Numeric = $true }
    @{ Cell = "E50"; Value = '  +0.97%  '; Numeric = $false }
    @{ Cell = "D51"; Value = '0.0845'; Numeric = $true }
    @{ Cell = "E51"; Value = '  -3.87%  '; Numeric = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.Numeric) {
        # Force text storage so the displayed digits (incl. trailing zeros)
        # are preserved exactly instead of being parsed into a Double.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
